# Daily attendance processing - 2025-11-01 19:42:10
# Normalize the "Recorded By" entries in column G: swap the order of the
# two known name/email pairs so the most recently-acting recorder is
# listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value2 = "System, backup@backdoor.com"
    }
}
